# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45179 to 45180 (i.e. 2023-09-10 -> 2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Data rows run from row 2 to row 270 (row 1 is the header row).
$firstRow = 2
$lastRow = 270

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2() -eq 45179) {
        $cell.Value = 45180
    }
}
